$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.326.21"
$ws.Range("E2").Value = "  +9.36%  "

$ws.Range("D3").Value = "3.227.66"
$ws.Range("E3").Value = "  +4.00%  "

$ws.Range("D5").Value = "'398.49"
$ws.Range("E5").Value = "  +4.18%  "

$ws.Range("D6").Value = "'110.92"
$ws.Range("E6").Value = "  +7.43%  "

$ws.Range("E7").Value = "  +2.19%  "

$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("D9").Value = "'0.627"
$ws.Range("E9").Value = "  +7.59%  "

$ws.Range("D10").Value = "'39.76"
$ws.Range("E10").Value = "  +6.73%  "

$ws.Range("E11").Value = "  +4.69%  "

$ws.Range("E12").Value = "  +2.08%  "

$ws.Range("D13").Value = "3.776.25"
$ws.Range("E13").Value = "  +4.96%  "

$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'19.13"
$ws.Range("E14").Value = "  +2.73%  "

$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "'8.07"
$ws.Range("E15").Value = "  +3.20%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.244.49"
$ws.Range("E16").Value = "  +4.71%  "

$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").Value = "'1.07"
$ws.Range("E17").Value = "  +7.28%  "

$ws.Range("D18").Value = "'10.48"
$ws.Range("E18").Value = "  -5.59%  "

$ws.Range("D19").Value = "56.065.54"
$ws.Range("E19").Value = "  +8.77%  "

$ws.Range("D20").Value = "'3.35"
$ws.Range("E20").Value = "  +2.13%  "

$ws.Range("D21").Value = "'13.11"
$ws.Range("E21").Value = "  +5.94%  "

$ws.Range("D22").Value = "'0.0000100"
$ws.Range("E22").Value = "  +4.30%  "

$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "'288.23"
$ws.Range("E23").Value = "  +8.23%  "

$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'73.98"
$ws.Range("E24").Value = "  +5.68%  "

$ws.Range("D25").Value = "'3.23"
$ws.Range("E25").Value = "  +4.85%  "

$ws.Range("D26").Value = "'8.18"
$ws.Range("E26").Value = "  +1.07%  "

$ws.Range("D27").Value = "'28.15"
$ws.Range("E27").Value = "  +4.22%  "

$ws.Range("E28").Value = "  +3.50%  "

$ws.Range("D29").Value = "'0.173"
$ws.Range("E29").Value = "  +2.61%  "

$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  -0.11%  "

$ws.Range("E31").Value = "  +4.42%  "

$ws.Range("D32").Value = "'11.27"
$ws.Range("E32").Value = "  +9.20%  "

$ws.Range("D33").Value = "'0.0501"
$ws.Range("E33").Value = "  +6.37%  "

$ws.Range("D34").Value = "'37.02"
$ws.Range("E34").Value = "  +4.39%  "

$ws.Range("E35").Value = "  +1.69%  "

$ws.Range("D36").Value = "'51.09"
$ws.Range("E36").Value = "  +1.66%  "

$ws.Range("D37").Value = "'3.60"
$ws.Range("E37").Value = "  +7.39%  "

$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  +0.12%  "

$ws.Range("D39").Value = "'3.08"
$ws.Range("E39").Value = "  +21.94%  "

$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").Value = "'139.54"
$ws.Range("E40").Value = "  +8.02%  "

$ws.Range("D41").Value = "'4.05"
$ws.Range("E41").Value = "  +10.71%  "

$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "'1.93"
$ws.Range("E42").Value = "  +2.30%  "

$ws.Range("D43").Value = "'0.284"
$ws.Range("E43").Value = "  -4.24%  "

$ws.Range("D44").Value = "'16.96"
$ws.Range("E44").Value = "  +1.33%  "

$ws.Range("E45").Value = "  +2.00%  "

$ws.Range("D46").Value = "'22.40"
$ws.Range("E46").Value = "  +0.33%  "

$ws.Range("E47").Value = "  +0.38%  "

$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("E48").Value = "  +0.94%  "

$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "2.125.33"
$ws.Range("E49").Value = "  +3.15%  "

$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").Value = "'2.07"
$ws.Range("E50").Value = "  +38.90%  "

$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "3.573.48"
$ws.Range("E51").Value = "  +4.60%  "
